# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on the
# zh-cn and de-de sheets for row 2 (the 528c1284-... handoff file).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 13:15:56"
$wsZhCn.Range("H2").Value = "2016-03-22 13:16:24"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 13:15:59"
$wsDeDe.Range("H2").Value = "2016-03-22 13:16:31"
